$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 34.20695657179154
$ws.Range("C2").Value = 16.48957345452475
$ws.Range("D2").Value = 6.084516429361718
$ws.Range("E2").Value = 7.400331956963808
$ws.Range("G2").Value = 3.785688023907315
$ws.Range("L2").Value = 11.36531711969576
$ws.Range("N2").Value = 24.44660392738255
$ws.Range("B3").Value = 33.65162825580455
$ws.Range("C3").Value = 15.87768574033236
$ws.Range("D3").Value = 5.985223694366822
$ws.Range("E3").Value = 7.383559142587797
$ws.Range("G3").Value = 3.792791826739354
$ws.Range("L3").Value = 11.35282454331879
$ws.Range("N3").Value = 24.34342644317945
$ws.Range("B4").Value = 33.31868454241758
$ws.Range("C4").Value = 15.49675428937713
$ws.Range("D4").Value = 5.925718959894797
$ws.Range("E4").Value = 7.373137838210761
$ws.Range("G4").Value = 3.797362093414316
$ws.Range("L4").Value = 11.34795006038904
$ws.Range("N4").Value = 24.28101900450909
$ws.Range("B5").Value = 33.1851921875388
$ws.Range("C5").Value = 15.34047527563503
$ws.Range("D5").Value = 5.901867252676617
$ws.Range("E5").Value = 7.368859395947845
$ws.Range("G5").Value = 3.799277252383894
$ws.Range("L5").Value = 11.34666524037174
$ws.Range("N5").Value = 24.25582803741583
$ws.Range("B6").Value = 33.16316259313255
$ws.Range("C6").Value = 15.31447036388212
$ws.Range("D6").Value = 5.897931513567617
$ws.Range("E6").Value = 7.368147026989134
$ws.Range("G6").Value = 3.799598457332781
$ws.Range("L6").Value = 11.34649421679427
$ws.Range("N6").Value = 24.25165971296945
$ws.Range("B7").Value = 33.31687515578983
$ws.Range("C7").Value = 15.49465052077567
$ws.Range("D7").Value = 5.925395641878235
$ws.Range("E7").Value = 7.37308026697329
$ws.Range("G7").Value = 3.797387707974131
$ws.Range("L7").Value = 11.34792989454595
$ws.Range("N7").Value = 24.28067829228753
$ws.Range("B8").Value = 34.01391344013457
$ws.Range("C8").Value = 16.27983417854958
$ws.Range("D8").Value = 6.04999490863867
$ws.Range("E8").Value = 7.39457366692749
$ws.Range("G8").Value = 3.788094326994687
$ws.Range("L8").Value = 11.3604285968195
$ws.Range("N8").Value = 24.41083442211718
$ws.Range("B9").Value = 35.43706442108305
$ws.Range("C9").Value = 17.76745708742209
$ws.Range("D9").Value = 6.304601568413521
$ws.Range("E9").Value = 7.435786749140624
$ws.Range("G9").Value = 3.771510270052011
$ws.Range("L9").Value = 11.40717237260056
$ws.Range("N9").Value = 24.6734128166783
$ws.Range("B10").Value = 36.50687486248272
$ws.Range("C10").Value = 18.81608946938698
$ws.Range("D10").Value = 6.496124126951035
$ws.Range("E10").Value = 7.465551343589564
$ws.Range("G10").Value = 3.760305974907577
$ws.Range("L10").Value = 11.45510955761834
$ws.Range("N10").Value = 24.87060066743078
$ws.Range("B11").Value = 36.99669994638118
$ws.Range("C11").Value = 19.2813262464533
$ws.Range("D11").Value = 6.583843683778287
$ws.Range("E11").Value = 7.478992581242434
$ws.Range("G11").Value = 3.755417317389702
$ws.Range("L11").Value = 11.47986803171416
$ws.Range("N11").Value = 24.96118470882153
$ws.Range("B12").Value = 37.18246634697305
$ws.Range("C12").Value = 19.4556504940798
$ws.Range("D12").Value = 6.61711556506634
$ws.Range("E12").Value = 7.484069189893313
$ws.Range("G12").Value = 3.753595708636907
$ws.Range("L12").Value = 11.48966682996629
$ws.Range("N12").Value = 24.99560843506217
$ws.Range("B13").Value = 37.14244848472683
$ws.Range("C13").Value = 19.41819122302595
$ws.Range("D13").Value = 6.609947937100258
$ws.Range("E13").Value = 7.482976427987436
$ws.Range("G13").Value = 3.753986712521203
$ws.Range("L13").Value = 11.48753768020559
$ws.Range("N13").Value = 24.98818935323353
$ws.Range("B14").Value = 37.01197842341667
$ws.Range("C14").Value = 19.29570591822162
$ws.Range("D14").Value = 6.586580066522471
$ws.Range("E14").Value = 7.479410498057021
$ws.Range("G14").Value = 3.755266860785278
$ws.Range("L14").Value = 11.48066570936962
$ws.Range("N14").Value = 24.96401436576222
$ws.Range("B15").Value = 36.93209315460219
$ws.Range("C15").Value = 19.22043483705476
$ws.Range("D15").Value = 6.572272745385447
$ws.Range("E15").Value = 7.477224558902607
$ws.Range("G15").Value = 3.756054836257197
$ws.Range("L15").Value = 11.47651151320943
$ws.Range("N15").Value = 24.94922214944
$ws.Range("B16").Value = 36.47491224737698
$ws.Range("C16").Value = 18.78543414537817
$ws.Range("D16").Value = 6.490400734699153
$ws.Range("E16").Value = 7.464671016257538
$ws.Range("G16").Value = 3.760629622737058
$ws.Range("L16").Value = 11.45355084093224
$ws.Range("N16").Value = 24.86469808927952
$ws.Range("B17").Value = 36.19513087226905
$ws.Range("C17").Value = 18.51543830921693
$ws.Range("D17").Value = 6.440304836163023
$ws.Range("E17").Value = 7.456945158096001
$ws.Range("G17").Value = 3.763489210026041
$ws.Range("L17").Value = 11.44022053544682
$ws.Range("N17").Value = 24.8130675647285
$ws.Range("B18").Value = 36.03451589683311
$ws.Range("C18").Value = 18.35904192456835
$ws.Range("D18").Value = 6.411548817741871
$ws.Range("E18").Value = 7.452492091127502
$ws.Range("G18").Value = 3.765153588400249
$ws.Range("L18").Value = 11.43283121493174
$ws.Range("N18").Value = 24.783454436786
$ws.Range("B19").Value = 35.98019270579373
$ws.Range("C19").Value = 18.30590460473363
$ws.Range("D19").Value = 6.401823418906366
$ws.Range("E19").Value = 7.450982724431507
$ws.Range("G19").Value = 3.765720497775037
$ws.Range("L19").Value = 11.43037707742547
$ws.Range("N19").Value = 24.77344242064441
$ws.Range("B20").Value = 36.22488349495936
$ws.Range("C20").Value = 18.54429506033777
$ws.Range("D20").Value = 6.445631874949917
$ws.Range("E20").Value = 7.457768550799935
$ws.Range("G20").Value = 3.763182773880241
$ws.Range("L20").Value = 11.44161081240184
$ws.Range("N20").Value = 24.81855513026773
$ws.Range("B21").Value = 37.05029440071635
$ws.Range("C21").Value = 19.33173418983542
$ws.Range("D21").Value = 6.593442544089245
$ws.Range("E21").Value = 7.480458251602842
$ws.Range("G21").Value = 3.754890048934076
$ws.Range("L21").Value = 11.48267269694558
$ws.Range("N21").Value = 24.97111189161463
$ws.Range("B22").Value = 37.59130490445597
$ws.Range("C22").Value = 19.83552501901535
$ws.Range("D22").Value = 6.690348650468738
$ws.Range("E22").Value = 7.49521021919954
$ws.Range("G22").Value = 3.749642767683335
$ws.Range("L22").Value = 11.51197547330912
$ws.Range("N22").Value = 25.07152218357219
$ws.Range("B23").Value = 37.30247117991408
$ws.Range("C23").Value = 19.56768096276641
$ws.Range("D23").Value = 6.63861032441439
$ws.Range("E23").Value = 7.487343560789404
$ws.Range("G23").Value = 3.752427665107643
$ws.Range("L23").Value = 11.4961108477137
$ws.Range("N23").Value = 25.01786858967576
$ws.Range("B24").Value = 36.21143158876535
$ws.Range("C24").Value = 18.53125256224559
$ws.Range("D24").Value = 6.443223378855797
$ws.Range("E24").Value = 7.457396330436511
$ws.Range("G24").Value = 3.763321250209176
$ws.Range("L24").Value = 11.44098141326037
$ws.Range("N24").Value = 24.81607398361044
$ws.Range("B25").Value = 35.04708234638405
$ws.Range("C25").Value = 17.37198205257585
$ws.Range("D25").Value = 6.234810002351569
$ws.Range("E25").Value = 7.424731960759559
$ws.Range("G25").Value = 3.775823147902824
$ws.Range("L25").Value = 11.39214054975137
$ws.Range("N25").Value = 24.60161215210496
